$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts the string into a
# numeric cell (changing both the stored type and the precision).
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D32", "D33", "D34", "D36", "D37", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.613.17'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '3.307.93'
$ws.Range('E3').Value = '  +1.66%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '604.88'
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').Value = '140.76'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.307.75'
$ws.Range('E8').Value = '  +1.82%  '
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').Value = '0.149'
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('D11').Value = '5.53'
$ws.Range('E11').Value = '  +2.98%  '
$ws.Range('D12').Value = '0.467'
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').Value = '0.0000245'
$ws.Range('E13').Value = '  -0.50%  '
$ws.Range('D14').Value = '34.66'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('D15').Value = '3.848.69'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').Value = '3.307.18'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '63.713.01'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '6.84'
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').Value = '479.32'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').Value = '13.98'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = '0.736'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').Value = '7.97'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').Value = '13.91'
$ws.Range('E24').Value = '  +5.74%  '
$ws.Range('D25').Value = '85.17'
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '2.77'
$ws.Range('E27').Value = '  +1.25%  '
$ws.Range('D29').Value = '7.18'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('D30').Value = '8.13'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('D32').Value = '28.64'
$ws.Range('E32').Value = '  +3.65%  '
$ws.Range('D33').Value = '0.105'
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('D34').Value = '2.51'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('D36').Value = '6.03'
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('D37').Value = '52.41'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('D38').Value = '0.0₃0734'
$ws.Range('E38').Value = '  +3.07%  '
$ws.Range('D39').Value = '0.0399'
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('D40').Value = '3.125.91'
$ws.Range('E40').Value = '  +4.76%  '
$ws.Range('D41').Value = '429.44'
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('D42').Value = '0.120'
$ws.Range('E42').Value = '  +8.58%  '
$ws.Range('D43').Value = '8.30'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').Value = '2.71'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').Value = '0.263'
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('D46').Value = '2.21'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('D47').Value = '36.53'
$ws.Range('E47').Value = '  +8.79%  '
$ws.Range('D48').Value = '26.25'
$ws.Range('E48').Value = '  +1.55%  '
$ws.Range('D49').Value = '127.63'
$ws.Range('E49').Value = '  +4.98%  '
$ws.Range('D51').Value = '2.30'
$ws.Range('E51').Value = '  -0.54%  '

# Restore the default (Normal) cell style so only the value changes,
# keeping the cell type as text without leaving a custom number format
# applied to the cell.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
